# Generate Report for Handoff
# Updates the localization-status report: Priority bumps from "low" to "ht"
# for the still-pending files, and the handoff timestamps for those same
# rows (zh-cn / de-de worksheets) plus the "Latest HO Xliff Generate Date"
# column on the Overview sheet are refreshed to the new handoff run time.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: rows 4-7 -> Priority (E) low -> ht ; Latest Handoff Datetime (H) refreshed
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-25 00:29:13"
}

# de-de: rows 4-7 -> Priority (E) low -> ht ; Latest Handoff Datetime (H) refreshed
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-25 00:29:18"
}

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G) refreshed
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-25 00:29:18"
}
